# Visit form workbook update:
#  - survey!F2:    "plot_name" -> "plot_id"
#  - queries!E2:   "plot_id >= ?" -> "_id >= ?"
#  - settings:     new row (table_id / visit)
#  - selection/active-sheet state updated to match the end of the edit session
#    (settings sheet ends up active, each sheet's last-used cell updated)

$wb = $excel.ActiveWorkbook

# --- survey sheet: rename the linked-plot field from plot_name to plot_id ---
$survey = $wb.Worksheets.Item("survey")
$survey.Range("F2").Value = "plot_id"
$survey.Range("F8").Select() | Out-Null

# --- choices sheet: untouched data, just revisit the previously selected cell ---
$choices = $wb.Worksheets.Item("choices")
$choices.Range("B18").Select() | Out-Null

# --- queries sheet: the linked-plot query now filters on the bare "_id" column ---
$queries = $wb.Worksheets.Item("queries")
$queries.Range("E2").Value = "_id >= ?"
$queries.Range("E3").Select() | Out-Null

# --- settings sheet: record the table_id setting (value = visit) for this form ---
$settings = $wb.Worksheets.Item("settings")
$settings.Range("A5").Value = "table_id"
$settings.Range("B5").Value = "visit"
$settings.Range("B6").Select() | Out-Null

Write-Host "done"
